# Update column F (dSF) values for specific rows per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 1
    11 = 3
    13 = 2
    17 = -3
    19 = 1
    20 = -4
    21 = 4
    23 = 10
    24 = -5
    25 = -1
    26 = -1
    28 = 6
    30 = -4
    32 = 2
    33 = 3
    34 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
